$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 11 (the footer row), shifting the
# footer row down to row 12. This also shifts the A11:E11 merge down to
# A12:E12 and extends the sheet dimension to A1:F12 automatically.
$ws.Rows.Item(11).Insert()

# Copy the formatting from row 10 (the last data row) onto the newly
# inserted row 11 so the new row matches the existing data rows' styles.
$ws.Range("A10:F10").Copy()
$ws.Range("A11:F11").PasteSpecial(-4122) # xlPasteFormats

# Update row 7 values.
$ws.Range("A7").Value = 1
$ws.Range("B7").Value = 91
$ws.Range("C7").Value = ""
$ws.Range("D7").Value = 2.14
$ws.Range("E7").Value = ""

# Update row 8 values.
$ws.Range("A8").Value = 10
$ws.Range("B8").Value = 81.8
$ws.Range("C8").Value = 21.87235698318771
$ws.Range("D8").Value = 2.3003
$ws.Range("E8").Value = 0.5982073312080948

# Update row 9 values.
$ws.Range("A9").Value = 3
$ws.Range("B9").Value = 131.6666666666667
$ws.Range("C9").Value = 37.52776749732568
$ws.Range("D9").Value = 2.755
$ws.Range("E9").Value = 0.1281600561797629

# Update row 10 values.
$ws.Range("A10").Value = 4
$ws.Range("B10").Value = 115.25
$ws.Range("C10").Value = 9.178779875342908
$ws.Range("D10").Value = 3.38875
$ws.Range("E10").Value = 0.1162163929916946

# Fill in the new row 11 values.
$ws.Range("A11").Value = 14
$ws.Range("B11").Value = 209.2142857142857
$ws.Range("C11").Value = 50.97688551827051
$ws.Range("D11").Value = 3.999214285714287
$ws.Range("E11").Value = 0.7594047444769265
